# Crunchbase_PrivateCompanyOverview_Template.pptx - slide 1 cleanup
#  - rename the placeholder "logo" rectangle and the stock chart to their
#    semantic automation names
#  - fill in the previously-empty "Headquarters" / "Investors" table cells
#    with their placeholder tokens

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Rename shapes -------------------------------------------------------
$s.Shapes.Item(3).Name = "companyLogo"   # was "Rectangle 3"
$s.Shapes.Item(4).Name = "stockChart"    # was "Chart 4"

# --- Populate the company-facts table ------------------------------------
$table = $s.Shapes.Item(6).Table

# Row 4 = "Headquarters:" -> fill in the HQ Location placeholder
$table.Cell(4, 2).Shape.TextFrame.TextRange.Text = "<HQ Location>"

# Row 5 = "Investors:" -> fill in the Investors placeholder
$table.Cell(5, 2).Shape.TextFrame.TextRange.Text = "<Investors>"
